$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.528.33"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "1.605.13"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("E6").Value = "  +6.97%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'26.99"
$ws.Range("D9").Value = "'43.43"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("D12").Value = "'0.0911"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").Value = "1.835.29"
$ws.Range("E13").Value = "  +3.24%  "
$ws.Range("D14").Value = "1.601.06"
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("D15").Value = "29.549.37"
$ws.Range("E15").Value = "  +3.39%  "
$ws.Range("E16").Value = "  +4.79%  "
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").Value = "'63.55"
$ws.Range("E18").Value = "  +3.86%  "
$ws.Range("D19").Value = "'242.98"
$ws.Range("E19").Value = "  +5.83%  "
$ws.Range("D20").Value = "'7.63"
$ws.Range("E20").Value = "  +3.49%  "
$ws.Range("E21").Value = "  +3.20%  "
$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("D24").Value = "'9.17"
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").Value = "'2.08"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "'154.76"
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("E27").Value = "  +3.87%  "
$ws.Range("E28").Value = "  +5.32%  "
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "'3.23"
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("D34").Value = "1.423.78"
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("E35").Value = "  +3.75%  "
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("E37").Value = "  +2.03%  "
$ws.Range("E38").Value = "  +5.91%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("E41").Value = "  +2.99%  "
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").Value = "'52.86"
$ws.Range("E43").Value = "  +20.89%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'0.794"
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("D47").Value = "'65.58"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").Value = "'5.29"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "1.747.27"
$ws.Range("E49").Value = "  +3.48%  "
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").Value = "'0.836"
$ws.Range("E51").Value = "  -3.57%  "
